# Apply the two textual edits found on slide 9 ("src/routes/api/category.js")
# of the Content Placeholder body.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(9)
$sh = $s.Shapes.Item(2)            # "Content Placeholder 8"
$tr = $sh.TextFrame.TextRange

# --- Edit 1 --------------------------------------------------------------
# "Other options could be role based bundling of the endpoint"
#   -> "...endpoints"  (paragraph 2, single run)
$para2 = $tr.Paragraphs(2)
$run2  = $para2.Runs(1)
$run2.Text = "Other options could be role based bundling of the endpoints"

# --- Edit 2 --------------------------------------------------------------
# Split the single run of paragraph 9 into three runs:
#   "Why aren’t all endpoints in this project with the new model? That would take time"
#   ", and the "
#   "five already demonstrate the new model well"
$para9 = $tr.Paragraphs(9)
$run9  = $para9.Runs(1)

$apostrophe = [char]0x2019
$run9.Text = "Why aren" + $apostrophe + "t all endpoints in this project with the new model? That would take time"

$midRun = $run9.InsertAfter(", and the ")
$endRun = $midRun.InsertAfter("five already demonstrate the new model well")
